# 2_Formulas_and_Functions_Solution.xlsx edit script
# - Remove the "Data" sheet
# - "1_Formulas_Intro": insert "2.1.3 - Total Compensation" (col J) and
#   "2.1.10 - Problem" (col M) formula columns
# - "2_Functions_Intro": insert "2.2.5 - Problem" (col K, row 2 only) and
#   fix the "2.2.3 Problem" / "2.1.7 Problem" header text (add the hyphen)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Remove the "Data" sheet entirely
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Data").Delete()

# ---------------------------------------------------------------------------
# 2. "1_Formulas_Intro" -- add Total Compensation + Problem 2.1.10 columns
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1_Formulas_Intro")

# Insert a new column at J; existing J (Meet Experience) and K (High Stock
# Options) shift right to K and L.
$ws1.Columns.Item(10).Insert()

# Header + formulas for the new "Total Compensation" column (J)
$ws1.Cells.Item(1, 10).Value = "2.1.3 - Total Compensation"
$ws1.Range("J2:J11").Formula = "=D2+E2+(D2*F2/100)"
$ws1.Range("J2:J11").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"
$ws1.Columns.Item(10).ColumnWidth = $ws1.Columns.Item(9).ColumnWidth

# Copy the header look (bold font + grey fill) from an existing header cell
$ws1.Cells.Item(1, 9).Copy()
$ws1.Cells.Item(1, 10).PasteSpecial(-4122)
$ws1.Cells.Item(1, 10).Value = "2.1.3 - Total Compensation"
$excel.CutCopyMode = 0

# New "2.1.10 - Problem" column (M) -- Boolean AND() formula
$ws1.Cells.Item(1, 13).Value = "2.1.10 - Problem"
$ws1.Range("M2:M11").Formula = "=AND(B2>=5,J2>=100000)"
$ws1.Columns.Item(13).ColumnWidth = 14.6

$ws1.Cells.Item(1, 11).Copy()
$ws1.Cells.Item(1, 13).PasteSpecial(-4122)
$ws1.Cells.Item(1, 13).Value = "2.1.10 - Problem"
$excel.CutCopyMode = 0

# Selection / active-tab bookkeeping (this sheet stays the active tab)
$ws1.Activate()
$ws1.Range("F23").Select()

# ---------------------------------------------------------------------------
# 3. "2_Functions_Intro" -- fix header text + add "2.2.5 - Problem" column
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2_Functions_Intro")

# Typo / formatting fixes on existing problem headers
$ws2.Cells.Item(1, 10).Value = "2.2.3 - Problem"
$ws2.Cells.Item(1, 11).Value = "2.1.7 - Problem"

# Insert a new column at K (pushes the old "2.1.7 Problem" COUNTIF column to L)
$ws2.Columns.Item(11).Insert()
$ws2.Cells.Item(1, 11).Value = "2.2.5 - Problem"
$ws2.Cells.Item(2, 11).Formula = "=MAX(H2:H11)"
$ws2.Columns.Item(11).ColumnWidth = $ws2.Columns.Item(10).ColumnWidth

$ws2.Cells.Item(1, 10).Copy()
$ws2.Cells.Item(1, 11).PasteSpecial(-4122)
$ws2.Cells.Item(1, 11).Value = "2.2.5 - Problem"
$excel.CutCopyMode = 0

$ws2.Cells.Item(2, 10).Copy()
$ws2.Cells.Item(2, 11).PasteSpecial(-4122)
$ws2.Cells.Item(2, 11).Formula = "=MAX(H2:H11)"
$excel.CutCopyMode = 0

$ws2.Range("J28").Select()

# Re-activate "1_Formulas_Intro" last so it remains the visible/active tab
$ws1.Activate()
$ws1.Range("F23").Select()

$wb.Application.Calculate()
